$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column (C) from 45174 to 45175
# for all data rows (2 through 26).
for ($row = 2; $row -le 26; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45174) {
        $cell.Value2 = 45175
    }
}
